$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): shift existing labels one column right, add "rownumber" ---
$ws.Range("F1").Value = "roletype"
$ws.Range("E1").Value = "rolecode"
$ws.Range("D1").Value = "email"
$ws.Range("C1").Value = "password"
$ws.Range("B1").Value = "username"
$ws.Range("A1").Value = "rownumber"

# --- Data row (row 2): shift existing values one column right, add row number ---
$ws.Range("F2").Value = "MANAGER"
$ws.Range("E2").Value = "MANAGER"
$ws.Range("D2").Value = "test@example.com"
$ws.Range("C2").Value = "pass11dasdsad"
$ws.Range("B2").Value = "testuser3"
$ws.Range("A2").Value = 1

# --- Header styling: bold, matching the rest of row 1 ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("F1").Font.Bold = $true

# --- Re-home the hyperlink from the old email cell (C2) to the new one (D2) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:test@example.com")
$ws.Range("D2").Style = "Hyperlink"
$ws.Range("C2").Style = "Normal"

# --- Column widths: A and D are new/custom; B/C keep their pre-existing widths ---
$ws.Range("A1").EntireColumn.ColumnWidth = 23
$ws.Range("D1").EntireColumn.ColumnWidth = 15.1666666667

# --- Selection shown in the saved view ---
$ws.Range("C13").Select()
